# The "busqueda" sheet used to list several scraped products (rows 3-11).
# Per the commit ("add en el objeto de product-url, product-img y add
# oncity en el scraping") the sample/test data produced by older scraping
# runs is cleaned out, leaving just the header row and a single sample
# row (Planchita Philips Hp8321) that was already in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete sample rows (3 through 11), shifting nothing up
# since they are the last rows of the used range.
$ws.Range("A3:E11").EntireRow.Delete() | Out-Null

# Reflect where the user was last positioned after trimming the sheet.
$ws.Range("A3").Select() | Out-Null
